$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 12588.059
$ws.Cells.Item(72, 8).Value = 12588.059
$ws.Cells.Item(111, 8).Value = 8771.75
$ws.Cells.Item(111, 9).Value = 5027.5
$ws.Cells.Item(111, 11).Value = 15082.5
$ws.Cells.Item(111, 13).Value = -12015.5
$ws.Cells.Item(138, 8).Value = 4673.8975
$ws.Cells.Item(138, 10).Value = 8554.6
$ws.Cells.Item(138, 12).Value = 25663.8
$ws.Cells.Item(138, 14).Value = -35943.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 5681024
$ws.Cells.Item(61, 9).Value = 6065966
$ws.Cells.Item(61, 11).Value = 6065966
$ws.Cells.Item(61, 13).Value = -6065754
$ws.Cells.Item(62, 8).Value = 97499.5
$ws.Cells.Item(62, 10).Value = 97499.5
$ws.Cells.Item(62, 12).Value = 97499.5
$ws.Cells.Item(62, 14).Value = -98747.5
$ws.Cells.Item(65, 8).Value = 97499.5
$ws.Cells.Item(65, 10).Value = 97499.5
$ws.Cells.Item(65, 12).Value = 292498.5
$ws.Cells.Item(65, 14).Value = -298738.5
$ws.Cells.Item(74, 8).Value = 1715.5358
$ws.Cells.Item(74, 9).Value = 1558.0435
$ws.Cells.Item(74, 10).Value = 2440
$ws.Cells.Item(74, 11).Value = 1558.0435
$ws.Cells.Item(74, 12).Value = 2440
$ws.Cells.Item(74, 13).Value = -684.0435
$ws.Cells.Item(74, 14).Value = -4188
$ws.Cells.Item(77, 8).Value = 1715.5358
$ws.Cells.Item(77, 9).Value = 1558.0435
$ws.Cells.Item(77, 10).Value = 2440
$ws.Cells.Item(77, 11).Value = 7790.2175
$ws.Cells.Item(77, 12).Value = 12200
$ws.Cells.Item(77, 13).Value = -3422.2175
$ws.Cells.Item(77, 14).Value = -20936
$ws.Cells.Item(94, 8).Value = 91748.75
$ws.Cells.Item(94, 10).Value = 91748.75
$ws.Cells.Item(94, 12).Value = 91748.75
$ws.Cells.Item(94, 14).Value = -93550.75
$ws.Cells.Item(97, 8).Value = 1203.8077
$ws.Cells.Item(97, 9).Value = 691.6316
$ws.Cells.Item(97, 11).Value = 691.6316
$ws.Cells.Item(97, 13).Value = -195.6316
$ws.Cells.Item(102, 8).Value = 2949.6924
$ws.Cells.Item(102, 9).Value = 2814.6
$ws.Cells.Item(102, 11).Value = 2814.6
$ws.Cells.Item(102, 13).Value = -1192.6
$ws.Cells.Item(110, 8).Value = 4988.1924
$ws.Cells.Item(110, 9).Value = 4857.263
$ws.Cells.Item(110, 11).Value = 4857.263
$ws.Cells.Item(110, 13).Value = -2812.263
$ws.Cells.Item(136, 8).Value = 5681024
$ws.Cells.Item(136, 9).Value = 6065966
$ws.Cells.Item(136, 11).Value = 18197898
$ws.Cells.Item(136, 13).Value = -18195348

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 7974.3335
$ws.Cells.Item(22, 9).Value = 1045
$ws.Cells.Item(22, 10).Value = 21833
$ws.Cells.Item(22, 11).Value = 1045
$ws.Cells.Item(22, 12).Value = 21833
$ws.Cells.Item(22, 13).Value = -872
$ws.Cells.Item(22, 14).Value = -22179
$ws.Cells.Item(82, 8).Value = 5249.5
$ws.Cells.Item(82, 9).Value = 5249.5
$ws.Cells.Item(82, 11).Value = 5249.5
$ws.Cells.Item(82, 13).Value = -4866.5
$ws.Cells.Item(85, 8).Value = 5249.5
$ws.Cells.Item(85, 9).Value = 5249.5
$ws.Cells.Item(85, 11).Value = 5249.5
$ws.Cells.Item(85, 13).Value = -3923.5
$ws.Cells.Item(99, 8).Value = 3099.875
$ws.Cells.Item(99, 10).Value = 3149.8333
$ws.Cells.Item(99, 12).Value = 3149.8333
$ws.Cells.Item(99, 14).Value = -6145.8333
$ws.Cells.Item(105, 8).Value = 448013.38
$ws.Cells.Item(105, 9).Value = 646124.9
$ws.Cells.Item(105, 11).Value = 646124.9
$ws.Cells.Item(105, 13).Value = -644377.9
$ws.Cells.Item(107, 8).Value = 4242.2144
$ws.Cells.Item(107, 9).Value = 4376.231
$ws.Cells.Item(107, 11).Value = 4376.231
$ws.Cells.Item(107, 13).Value = -2456.231

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 4765.5557
$ws.Cells.Item(58, 10).Value = 6503.25
$ws.Cells.Item(58, 12).Value = 6503.25
$ws.Cells.Item(58, 14).Value = -6909.25
$ws.Cells.Item(63, 8).Value = 79507.75
$ws.Cells.Item(63, 10).Value = 79507.75
$ws.Cells.Item(63, 12).Value = 79507.75
$ws.Cells.Item(63, 14).Value = -80879.75
$ws.Cells.Item(66, 8).Value = 79507.75
$ws.Cells.Item(66, 10).Value = 79507.75
$ws.Cells.Item(66, 12).Value = 238523.25
$ws.Cells.Item(66, 14).Value = -245387.25
$ws.Cells.Item(105, 8).Value = 5219.1665
$ws.Cells.Item(105, 9).Value = 1072
$ws.Cells.Item(105, 11).Value = 1072
$ws.Cells.Item(105, 13).Value = 675
$ws.Cells.Item(107, 8).Value = 2057.1428
$ws.Cells.Item(107, 9).Value = 832.4286
$ws.Cells.Item(107, 11).Value = 832.4286
$ws.Cells.Item(107, 13).Value = 1087.5714
$ws.Cells.Item(132, 8).Value = 2599.6428
$ws.Cells.Item(132, 9).Value = 2684.3076
$ws.Cells.Item(132, 10).Value = 1499
$ws.Cells.Item(132, 11).Value = 8052.9228
$ws.Cells.Item(132, 12).Value = 4497
$ws.Cells.Item(132, 13).Value = -5522.9228
$ws.Cells.Item(132, 14).Value = -9557
$ws.Cells.Item(134, 8).Value = 2830.1875
$ws.Cells.Item(134, 9).Value = 2607.0833
$ws.Cells.Item(134, 11).Value = 7821.249899999999
$ws.Cells.Item(134, 13).Value = -5286.249899999999
$ws.Cells.Item(136, 8).Value = 4765.5557
$ws.Cells.Item(136, 10).Value = 6503.25
$ws.Cells.Item(136, 12).Value = 19509.75
$ws.Cells.Item(136, 14).Value = -24609.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1007
$ws.Cells.Item(5, 9).Value = 834.75
$ws.Cells.Item(5, 10).Value = 1121.8334
$ws.Cells.Item(5, 11).Value = 2504.25
$ws.Cells.Item(5, 12).Value = 3365.5002
$ws.Cells.Item(5, 13).Value = -2392.25
$ws.Cells.Item(5, 14).Value = -3589.5002
$ws.Cells.Item(8, 8).Value = 245.2
$ws.Cells.Item(8, 9).Value = 245.2
$ws.Cells.Item(8, 11).Value = 735.5999999999999
$ws.Cells.Item(8, 13).Value = -596.5999999999999
$ws.Cells.Item(42, 8).Value = 1000000000
$ws.Cells.Item(42, 9).Value = 1000000000
$ws.Cells.Item(42, 11).Value = 3000000000
$ws.Cells.Item(42, 13).Value = -2999999466
$ws.Cells.Item(114, 8).Value = 3874.7856
$ws.Cells.Item(114, 10).Value = 10468
$ws.Cells.Item(114, 12).Value = 31404
$ws.Cells.Item(114, 14).Value = -37912
$ws.Cells.Item(135, 8).Value = 1007
$ws.Cells.Item(135, 9).Value = 834.75
$ws.Cells.Item(135, 10).Value = 1121.8334
$ws.Cells.Item(135, 11).Value = 7512.75
$ws.Cells.Item(135, 12).Value = 10096.5006
$ws.Cells.Item(135, 13).Value = -4977.75
$ws.Cells.Item(135, 14).Value = -15166.5006

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1800.5238
$ws.Cells.Item(97, 9).Value = 1745.25
$ws.Cells.Item(97, 10).Value = 1977.4
$ws.Cells.Item(97, 11).Value = 1745.25
$ws.Cells.Item(97, 12).Value = 1977.4
$ws.Cells.Item(97, 13).Value = -1249.25
$ws.Cells.Item(97, 14).Value = -2969.4
$ws.Cells.Item(113, 8).Value = 3708593.5
$ws.Cells.Item(113, 9).Value = 1169.5
$ws.Cells.Item(113, 10).Value = 6180209.5
$ws.Cells.Item(113, 11).Value = 1169.5
$ws.Cells.Item(113, 12).Value = 6180209.5
$ws.Cells.Item(113, 13).Value = 1000.5
$ws.Cells.Item(113, 14).Value = -6184549.5
$ws.Cells.Item(126, 8).Value = 4794.241
$ws.Cells.Item(126, 9).Value = 4359.9473
$ws.Cells.Item(126, 10).Value = 5619.4
$ws.Cells.Item(126, 11).Value = 13079.8419
$ws.Cells.Item(126, 12).Value = 16858.2
$ws.Cells.Item(126, 13).Value = -10609.8419
$ws.Cells.Item(126, 14).Value = -21798.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 5723.4443
$ws.Cells.Item(40, 9).Value = 4624.4
$ws.Cells.Item(40, 10).Value = 8863.571
$ws.Cells.Item(40, 11).Value = 4624.4
$ws.Cells.Item(40, 12).Value = 8863.571
$ws.Cells.Item(40, 13).Value = -4488.4
$ws.Cells.Item(40, 14).Value = -9135.571
$ws.Cells.Item(46, 8).Value = 940.625
$ws.Cells.Item(46, 9).Value = 696.3077
$ws.Cells.Item(46, 11).Value = 696.3077
$ws.Cells.Item(46, 13).Value = -508.3077
$ws.Cells.Item(122, 8).Value = 4114.909
$ws.Cells.Item(122, 9).Value = 3290.5557
$ws.Cells.Item(122, 10).Value = 7824.5
$ws.Cells.Item(122, 11).Value = 9871.667099999999
$ws.Cells.Item(122, 12).Value = 23473.5
$ws.Cells.Item(122, 13).Value = -7421.667099999999
$ws.Cells.Item(122, 14).Value = -28373.5
$ws.Cells.Item(136, 8).Value = 5926.077
$ws.Cells.Item(136, 9).Value = 3109.6
$ws.Cells.Item(136, 10).Value = 7686.375
$ws.Cells.Item(136, 11).Value = 9328.799999999999
$ws.Cells.Item(136, 12).Value = 23059.125
$ws.Cells.Item(136, 13).Value = -6778.799999999999
$ws.Cells.Item(136, 14).Value = -28159.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 6682.2666
$ws.Cells.Item(96, 10).Value = 5342.5713
$ws.Cells.Item(96, 12).Value = 5342.5713
$ws.Cells.Item(96, 14).Value = -8088.5713
$ws.Cells.Item(132, 8).Value = 1252420.1
$ws.Cells.Item(132, 9).Value = 2466.5454
$ws.Cells.Item(132, 11).Value = 7399.6362
$ws.Cells.Item(132, 13).Value = -4869.6362
$ws.Cells.Item(136, 8).Value = 272565.8
$ws.Cells.Item(136, 9).Value = 2496.577
$ws.Cells.Item(136, 11).Value = 7489.731000000001
$ws.Cells.Item(136, 13).Value = -4939.731000000001
